# "Generate Report for Handoff" - update the b.md row across all three
# sheets (Overview, zh-cn, de-de) now that the handback report has been
# generated for b.md, and widen the "Error Detail" column to show the
# new warning message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/73ecf6d5149dc312eb14272acdd06f9bce8628f5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/587fe92834fab02ad178f31ee371c3876260c814/e2e/b.md."

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------
# Overview sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-11-29 04:18:20"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-11-29 04:18:05"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 40 - (5/6)

# ---------------------------------------------------------------
# de-de sheet: row 3 corresponds to b.md
# ---------------------------------------------------------------
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-11-29 04:18:20"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 40 - (5/6)
